$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff Datetime (D4) and Correspond Handback DateTime (G4)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-19 05:38:51"
$wsZhCn.Range("G4").Value = "2016-02-19 05:39:34"

# de-de sheet: update Correspond Handoff Datetime (D4) and Correspond Handback DateTime (G4)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-19 05:39:01"
$wsDeDe.Range("G4").Value = "2016-02-19 05:39:50"
